$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value2 = 3
$ws.Range("G2").Value2 = 11.949039
$ws.Range("H2").Value2 = 35.847117
$ws.Range("I2").Value2 = 0.3580329233804654
$ws.Range("J2").Value2 = 0.3580329233804654
$ws.Range("K2").Value2 = 3
$ws.Range("M2").Value2 = 1.456609
$ws.Range("N2").Value2 = 4.369827
$ws.Range("O2").Value2 = 0.3001553107647443
$ws.Range("P2").Value2 = 0.3001553107647442
$ws.Range("Q2").Value2 = 17.405077748751
$ws.Range("R2").Value2 = 156.645699738759
$ws.Range("S2").Value2 = 0.1074654833812735
$ws.Range("T2").Value2 = 0.1074654833812734
$ws.Range("E3").Value2 = 3
$ws.Range("G3").Value2 = 11.949039
$ws.Range("H3").Value2 = 35.847117
$ws.Range("I3").Value2 = 0.3580329233804654
$ws.Range("J3").Value2 = 0.3580329233804654
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 2.339991
$ws.Range("N3").Value2 = 7.019973
$ws.Range("O3").Value2 = 0.482188923583271
$ws.Range("P3").Value2 = 0.4821889235832709
$ws.Range("Q3").Value2 = 27.960643718649
$ws.Range("R3").Value2 = 251.645793467841
$ws.Range("S3").Value2 = 0.1726395099321983
$ws.Range("T3").Value2 = 0.1726395099321983
$ws.Range("E4").Value2 = 3
$ws.Range("G4").Value2 = 11.949039
$ws.Range("H4").Value2 = 35.847117
$ws.Range("I4").Value2 = 0.3580329233804654
$ws.Range("J4").Value2 = 0.3580329233804654
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 1.056251
$ws.Range("N4").Value2 = 3.168753
$ws.Range("O4").Value2 = 0.2176557656519848
$ws.Range("P4").Value2 = 0.2176557656519848
$ws.Range("Q4").Value2 = 12.621184392789
$ws.Range("R4").Value2 = 113.590659535101
$ws.Range("S4").Value2 = 0.07792793006699361
$ws.Range("T4").Value2 = 0.0779279300669936
$ws.Range("E5").Value2 = 3
$ws.Range("G5").Value2 = 13.30334766666667
$ws.Range("H5").Value2 = 39.910043
$ws.Range("I5").Value2 = 0.3986125123403949
$ws.Range("J5").Value2 = 0.3986125123403949
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 1.456609
$ws.Range("N5").Value2 = 4.369827
$ws.Range("O5").Value2 = 0.3001553107647443
$ws.Range("P5").Value2 = 0.3001553107647442
$ws.Range("Q5").Value2 = 19.37777594139567
$ws.Range("R5").Value2 = 174.399983472561
$ws.Range("S5").Value2 = 0.1196456625162467
$ws.Range("T5").Value2 = 0.1196456625162467
$ws.Range("E6").Value2 = 3
$ws.Range("G6").Value2 = 13.30334766666667
$ws.Range("H6").Value2 = 39.910043
$ws.Range("I6").Value2 = 0.3986125123403949
$ws.Range("J6").Value2 = 0.3986125123403949
$ws.Range("K6").Value2 = 3
$ws.Range("M6").Value2 = 2.339991
$ws.Range("N6").Value2 = 7.019973
$ws.Range("O6").Value2 = 0.482188923583271
$ws.Range("P6").Value2 = 0.4821889235832709
$ws.Range("Q6").Value2 = 31.129713809871
$ws.Range("R6").Value2 = 280.167424288839
$ws.Range("S6").Value2 = 0.1922065382522384
$ws.Range("T6").Value2 = 0.1922065382522383
$ws.Range("E7").Value2 = 3
$ws.Range("G7").Value2 = 13.30334766666667
$ws.Range("H7").Value2 = 39.910043
$ws.Range("I7").Value2 = 0.3986125123403949
$ws.Range("J7").Value2 = 0.3986125123403949
$ws.Range("K7").Value2 = 3
$ws.Range("M7").Value2 = 1.056251
$ws.Range("N7").Value2 = 3.168753
$ws.Range("O7").Value2 = 0.2176557656519848
$ws.Range("P7").Value2 = 0.2176557656519848
$ws.Range("Q7").Value2 = 14.05167427626433
$ws.Range("R7").Value2 = 126.465068486379
$ws.Range("S7").Value2 = 0.0867603115719099
$ws.Range("T7").Value2 = 0.08676031157190989
$ws.Range("E8").Value2 = 3
$ws.Range("G8").Value2 = 8.121747999999998
$ws.Range("H8").Value2 = 24.365244
$ws.Range("I8").Value2 = 0.2433545642791397
$ws.Range("J8").Value2 = 0.2433545642791398
$ws.Range("K8").Value2 = 3
$ws.Range("M8").Value2 = 1.456609
$ws.Range("N8").Value2 = 4.369827
$ws.Range("O8").Value2 = 0.3001553107647443
$ws.Range("P8").Value2 = 0.3001553107647442
$ws.Range("Q8").Value2 = 11.830211232532
$ws.Range("R8").Value2 = 106.471901092788
$ws.Range("S8").Value2 = 0.07304416486722413
$ws.Range("T8").Value2 = 0.07304416486722412
$ws.Range("E9").Value2 = 3
$ws.Range("G9").Value2 = 8.121747999999998
$ws.Range("H9").Value2 = 24.365244
$ws.Range("I9").Value2 = 0.2433545642791397
$ws.Range("J9").Value2 = 0.2433545642791398
$ws.Range("K9").Value2 = 3
$ws.Range("M9").Value2 = 2.339991
$ws.Range("N9").Value2 = 7.019973
$ws.Range("O9").Value2 = 0.482188923583271
$ws.Range("P9").Value2 = 0.4821889235832709
$ws.Range("Q9").Value2 = 19.004817224268
$ws.Range("R9").Value2 = 171.043355018412
$ws.Range("S9").Value2 = 0.1173428753988343
$ws.Range("T9").Value2 = 0.1173428753988343
$ws.Range("E10").Value2 = 3
$ws.Range("G10").Value2 = 8.121747999999998
$ws.Range("H10").Value2 = 24.365244
$ws.Range("I10").Value2 = 0.2433545642791397
$ws.Range("J10").Value2 = 0.2433545642791398
$ws.Range("K10").Value2 = 3
$ws.Range("M10").Value2 = 1.056251
$ws.Range("N10").Value2 = 3.168753
$ws.Range("O10").Value2 = 0.2176557656519848
$ws.Range("P10").Value2 = 0.2176557656519848
$ws.Range("Q10").Value2 = 8.578604446747999
$ws.Range("R10").Value2 = 77.20744002073199
$ws.Range("S10").Value2 = 0.05296752401308132
$ws.Range("T10").Value2 = 0.05296752401308131
